$wb = $excel.ActiveWorkbook

# --- 1) "Horas Extra" sheet: insert a new "periodo_pago" column (D) ---
$wsHoras = $wb.Worksheets.Item("Horas Extra")
$wsHoras.Columns("D:D").Insert()
$wsHoras.Range("D1").Value = "periodo_pago"
$wsHoras.Range("D1").HorizontalAlignment = 1

# --- 2) "Asignaciones" sheet: insert a new "periodo_pago" column (D) ---
$wsAsig = $wb.Worksheets.Item("Asignaciones")
$wsAsig.Columns("D:D").Insert()
$wsAsig.Range("D1").Value = "periodo_pago"
$wsAsig.Range("D1").HorizontalAlignment = 1

# --- 3) "Guía" sheet: document the two new fields ("periodo" and "periodo_pago") ---
$wsGuia = $wb.Worksheets.Item("Guía")
$wsGuia.Rows("14:15").Insert()
$wsGuia.Range("A14").Value = "periodo"
$wsGuia.Range("B14").Value = "Fecha efectiva (día real de la asignacion)"
$wsGuia.Range("A15").Value = "periodo_pago"
$wsGuia.Range("B15").Value = "Fecha de corte/quincena que paga la asignación"
$wsGuia.Range("B14:B15").VerticalAlignment = -4108

# --- 4) "Horas Extra" sheet: mark the remaining header cells as text-formatted ---
$wsHoras.Range("E1:J1").NumberFormat = "@"

# --- selections / active sheet ---
$wsGuia.Range("B17").Select()
$wsAsig.Range("D10").Select()
$wsHoras.Range("F5").Select()
$wsHoras.Activate()
